# MSP430 Todo list edits
# 1. Split "MSP" + " 430 To do list" into "MSP" + " 430 " + [gramStart]"To"[gramEnd] + " do list"
# 2. Mark "Resitor" with a spellStart/spellEnd proofErr pair
# 3. Insert three new tabbed bullet items before "Do PCB layout and routing",
#    moving the bookmark "_GoBack" onto the new last inserted item, and mark
#    "Digikey"/"Samtec" with spellStart/spellEnd proofErr pairs.

$d = $word.ActiveDocument

# --- Edit 1: paragraph 1 ("MSP 430 To do list") ---
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range

$xml1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t>MSP</w:t></w:r>
<w:r><w:t xml:space="preserve"> 430 </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>To</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t xml:space="preserve"> do list</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$r1.InsertXML($xml1)

# --- Edit 2: paragraph 3 ("Resitor footprint") ---
$p3 = $d.Paragraphs.Item(3)
$r3 = $p3.Range

$xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>Resitor</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> footprint</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$r3.InsertXML($xml2)

# --- Edit 3: reorganize the tail (paragraphs "Do PCB layout and routing" through
#     "Submit Samtec order and finalize") ---
$pStart = $d.Paragraphs.Item(8)
$pEnd = $d.Paragraphs.Item(11)
$r = $d.Range($pStart.Range.Start, $pEnd.Range.End)

$xml3 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:tab/><w:t>-Add power and ground for MSP430</w:t></w:r></w:p>
<w:p><w:r><w:tab/><w:t xml:space="preserve"> -I would collect the powers and grounds together in the schematic diagram</w:t></w:r></w:p>
<w:p><w:r><w:tab/><w:t>-Fix the 3.3 label to 5V coming off of the USB</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
<w:p><w:r><w:t>Do PCB layout and routing</w:t></w:r></w:p>
<w:p><w:r><w:t>Complete Advance Circuits PCB check thingy</w:t></w:r></w:p>
<w:p><w:r><w:t xml:space="preserve">Submit order to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Digikey</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and finalize</w:t></w:r></w:p>
<w:p><w:r><w:t xml:space="preserve">Submit </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Samtec</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> order and finalize</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$r.InsertXML($xml3)

Write-Output "MSP430 todo list updated"
